# #12 Web interface schemas correctly placed on slides on presentation
#
# Repositions the four screenshot pictures on slide 18 ("Result") so the
# web-interface schema images line up correctly on the slide.
#
# NOTE on precision: PowerPoint's Shape.Left/Top/Width/Height are COM
# `Single` (32-bit float) properties expressed in points, while the
# underlying OOXML stores integer EMU (1 pt = 12700 EMU). Assigning a
# plain "EMU / 12700" double is not always enough to land back on the
# exact target EMU once the assigned value is rounded to a 32-bit float
# and re-quantized, so the literals below were chosen (from the middle of
# the float32 range that maps back to the exact target EMU) to reproduce
# the precise EMU coordinates from the target OOXML.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)

# "Image 3" (id 4) - top-left schema screenshot: moved right, ext unchanged
$imgResult = $s.Shapes.Item(3)
$imgResult.Left = 355.1367492675781
$imgResult.Top  = 21.046890258789062

# "Image 5" (id 6) - small schema screenshot: moved to top-left row, resized
$imgWeb1 = $s.Shapes.Item(4)
$imgWeb1.Left   = 0.0
$imgWeb1.Top    = 204.0886993408203
$imgWeb1.Width  = 304.5395812988281
$imgWeb1.Height = 129.4893341064453

# "Image 7" (id 8) - schema screenshot: moved to middle of row, resized
$imgWeb2 = $s.Shapes.Item(5)
$imgWeb2.Left   = 321.0926513671875
$imgWeb2.Top    = 204.0886993408203
$imgWeb2.Width  = 304.5395812988281
$imgWeb2.Height = 129.4893341064453

# "Image 9" (id 10) - schema screenshot: moved to right of row, resized
$imgWeb3 = $s.Shapes.Item(6)
$imgWeb3.Left   = 641.4595947265625
$imgWeb3.Top    = 204.0886993408203
$imgWeb3.Width  = 304.5396423339844
$imgWeb3.Height = 129.4893341064453
